# Update column F (dSF) values on Sheet1 to reflect repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -9
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = -3
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = -2
